$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2941, 14, 30, 33, 35, 48, 51),
    @(2942, 12, 30, 40, 46, 54, 60),
    @(2943, 8, 29, 30, 36, 39, 60),
    @(2944, 8, 15, 23, 39, 40, 59),
    @(2945, 1, 2, 3, 7, 27, 33)
)

$startRow = 398
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($c = 0; $c -lt 7; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $data[$i][$c]
    }
}

$ws.Activate() | Out-Null
$ws.Range("B398:G402").Select() | Out-Null
